$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text number format on cells whose new values look like
# numbers, so Excel keeps them as literal text (matches source formatting,
# e.g. trailing zeros like "8.00" or small decimals like "0.0000109").
$textCells = @("D5", "D6", "D7", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D32", "D33", "D34", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D44", "D45", "D47", "D48", "D49", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values row by row.
$ws.Range("D2").Value = '43.020.47'
$ws.Range("E2").Value = '  +1.20%  '

$ws.Range("D3").Value = '2.382.50'
$ws.Range("E3").Value = '  +6.60%  '

$ws.Range("E4").Value = '  -0.33%  '

$ws.Range("D5").Value = '329.52'
$ws.Range("E5").Value = '  +12.04%  '

$ws.Range("D6").Value = '105.74'
$ws.Range("E6").Value = '  -5.54%  '

$ws.Range("D7").Value = '0.642'
$ws.Range("E7").Value = '  +3.40%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  +10.00%  '

$ws.Range("D10").Value = '41.58'
$ws.Range("E10").Value = '  -4.16%  '

$ws.Range("D11").Value = '0.0939'
$ws.Range("E11").Value = '  +2.10%  '

$ws.Range("D12").Value = '8.63'
$ws.Range("E12").Value = '  -1.75%  '

$ws.Range("D13").Value = '1.04'
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").Value = '17.12'
$ws.Range("E14").Value = '  +14.38%  '

$ws.Range("E15").Value = '  +2.24%  '

$ws.Range("D16").Value = '2.744.29'

$ws.Range("D17").Value = '2.387.26'
$ws.Range("E17").Value = '  +6.65%  '

$ws.Range("D18").Value = '43.148.90'
$ws.Range("E18").Value = '  +1.57%  '

$ws.Range("D19").Value = '8.00'
$ws.Range("E19").Value = '  +11.43%  '

$ws.Range("D20").Value = '0.0000109'
$ws.Range("E20").Value = '  +2.80%  '

$ws.Range("D21").Value = '76.77'
$ws.Range("E21").Value = '  +3.68%  '

$ws.Range("D22").Value = '3.70'
$ws.Range("E22").Value = '  +10.69%  '

$ws.Range("D23").Value = '278.69'
$ws.Range("E23").Value = '  +14.01%  '

$ws.Range("D24").Value = '2.43'
$ws.Range("E24").Value = '  +2.77%  '

$ws.Range("D25").Value = '9.61'
$ws.Range("E25").Value = '  +7.90%  '

$ws.Range("D26").Value = '11.77'
$ws.Range("E26").Value = '  +3.23%  '

$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.08%  '

$ws.Range("D28").Value = '23.18'
$ws.Range("E28").Value = '  +7.27%  '

$ws.Range("D29").Value = '38.04'
$ws.Range("E29").Value = '  +2.55%  '

$ws.Range("D30").Value = '174.87'
$ws.Range("E30").Value = '  -0.21%  '

$ws.Range("E31").Value = '  -1.46%  '

$ws.Range("D32").Value = '3.16'
$ws.Range("E32").Value = '  +1.43%  '

$ws.Range("D33").Value = '0.0927'
$ws.Range("E33").Value = '  +5.49%  '

$ws.Range("D34").Value = '5.87'
$ws.Range("E34").Value = '  +3.22%  '

$ws.Range("E35").Value = '  +5.88%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '4.33'
$ws.Range("E36").Value = '  +2.57%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '4.88'
$ws.Range("E37").Value = '  -2.27%  '

$ws.Range("D38").Value = '0.0369'
$ws.Range("E38").Value = '  -0.93%  '

$ws.Range("D39").Value = '0.107'
$ws.Range("E39").Value = '  +2.87%  '

$ws.Range("D40").Value = '2.83'
$ws.Range("E40").Value = '  +17.99%  '

$ws.Range("E41").Value = '  +21.34%  '

$ws.Range("D42").Value = '0.235'
$ws.Range("E42").Value = '  +2.96%  '

$ws.Range("D43").Value = '70.03'
$ws.Range("E43").Value = '  -1.31%  '

$ws.Range("D44").Value = '123.77'
$ws.Range("E44").Value = '  +21.13%  '

$ws.Range("D45").Value = '94.48'
$ws.Range("E45").Value = '  +63.39%  '

$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D47").Value = '12.54'
$ws.Range("E47").Value = '  +1.67%  '

$ws.Range("D48").Value = '9.43'
$ws.Range("E48").Value = '  +10.78%  '

$ws.Range("D49").Value = '5.54'
$ws.Range("E49").Value = '  +1.89%  '

$ws.Range("D50").Value = '1.32'
$ws.Range("E50").Value = '  +2.30%  '

$ws.Range("D51").Value = '1.598.59'
$ws.Range("E51").Value = '  +11.86%  '

